$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - full row update
$ws.Range("A2").Value = "'" + " 47"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 154
$ws.Range("C2").Value = 0.96
$ws.Range("D2").Value = 181.7
$ws.Range("E2").Value = 0.83
$ws.Range("F2").Value = 126.1
$ws.Range("G2").Value = 0.74
$ws.Range("H2").Value = 0.62
$ws.Range("I2").Value = 0.4
$ws.Range("J2").Value = 0.08
$ws.Range("K2").Value = 0.12
$ws.Range("L2").Value = 95
$ws.Range("M2").Value = 114
$ws.Range("N2").Value = 62
$ws.Range("O2").Value = 13
$ws.Range("P2").Value = 19

# Row 3 - only A changes
$ws.Range("A3").Value = "'" + " 27"
$ws.Range("A3").Style = "Normal"

# Row 4 - only A changes
$ws.Range("A4").Value = "'" + " 13"
$ws.Range("A4").Style = "Normal"

# Row 5 - full row update
$ws.Range("A5").Value = "'" + " 13"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 48
$ws.Range("C5").Value = 0.88
$ws.Range("D5").Value = 153
$ws.Range("E5").Value = 0.8100000000000001
$ws.Range("F5").Value = 102.4
$ws.Range("G5").Value = 0.75
$ws.Range("H5").Value = 0.54
$ws.Range("I5").Value = 0.33
$ws.Range("J5").Value = 0.1
$ws.Range("K5").Value = 0.04
$ws.Range("L5").Value = 26
$ws.Range("M5").Value = 32
$ws.Range("N5").Value = 16
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 2
